# Insert two new data rows (a "Primera"/"Segunda" quality pair for Acelga)
# right before the current row 442, pushing the existing rows 442..558 down
# to 444..560 (matches the new dimension A1:R560).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("442:443").Insert()

# New row 442 - Calidad "Primera"
$ws.Range("A442").Value = 8
$ws.Range("B442").Value = "Terminal La Palmera de La Serena"
$ws.Range("C442").Value = "Coquimbo"
$ws.Range("D442").Value = 44932
$ws.Range("E442").Value = 4
$ws.Range("F442").Value = 100112009
$ws.Range("G442").Value = "Acelga"
$ws.Range("H442").Value = "Sin especificar"
$ws.Range("I442").Value = "Primera"
$ws.Range("J442").Value = 2300
$ws.Range("K442").Value = 600
$ws.Range("L442").Value = 700
$ws.Range("M442").Value = 650
$ws.Range("N442").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O442").Value = "Provincia del Elquí"
$ws.Range("P442").Value = 325
$ws.Range("Q442").Value = 2
$ws.Range("R442").Value = "Hortaliza"

# New row 443 - Calidad "Segunda"
$ws.Range("A443").Value = 8
$ws.Range("B443").Value = "Terminal La Palmera de La Serena"
$ws.Range("C443").Value = "Coquimbo"
$ws.Range("D443").Value = 44932
$ws.Range("E443").Value = 4
$ws.Range("F443").Value = 100112009
$ws.Range("G443").Value = "Acelga"
$ws.Range("H443").Value = "Sin especificar"
$ws.Range("I443").Value = "Segunda"
$ws.Range("J443").Value = 1540
$ws.Range("K443").Value = 500
$ws.Range("L443").Value = 550
$ws.Range("M443").Value = 525
$ws.Range("N443").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O443").Value = "Provincia del Elquí"
$ws.Range("P443").Value = 262
$ws.Range("Q443").Value = 2
$ws.Range("R443").Value = "Hortaliza"
